$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Labor Productivity 2012/2015 (columns R and S): divide the raw peso
# values by 100 by replacing the literal numbers with "=<value>/100"
# formulas, for every data row (2 through 18). ---
for ($i = 2; $i -le 18; $i++) {
    $rCell = $ws.Range("R$i")
    $rRaw = $rCell.Value2
    $rCell.Formula = "=$rRaw/100"

    $sCell = $ws.Range("S$i")
    $sRaw = $sCell.Value2
    $sCell.Formula = "=$sRaw/100"
}

# --- Center the "AverageDailyWage"-type column (G2:G18), which previously
# had no horizontal alignment applied. ---
$ws.Range("G2:G18").HorizontalAlignment = -4108

# --- Update the frozen-pane scroll position / active selection on the
# sheet view: scroll the right-hand (unfrozen) pane so column N is the
# left-most visible column, and move the active cell/selection to S13. ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 14
$ws.Range("S13").Select()
